$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 398.75
$ws.Range("J33").Value = 431.66666
$ws.Range("L33").Value = 431.66666
$ws.Range("N33").Value = -889.66666
$ws.Range("H94").Value = 16295.2
$ws.Range("I94").Value = 9772.444
$ws.Range("J94").Value = 75000
$ws.Range("K94").Value = 9772.444
$ws.Range("L94").Value = 75000
$ws.Range("M94").Value = -9321.444
$ws.Range("N94").Value = -75902
$ws.Range("H132").Value = 1531.421
$ws.Range("I132").Value = 1531.421
$ws.Range("K132").Value = 4594.263
$ws.Range("M132").Value = -2064.263

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 20200
$ws.Range("J10").Value = 20333.334
$ws.Range("L10").Value = 20333.334
$ws.Range("N10").Value = -20673.334
$ws.Range("H32").Value = 8034.35
$ws.Range("I32").Value = 5811.5757
$ws.Range("J32").Value = 18513.143
$ws.Range("K32").Value = 5811.5757
$ws.Range("L32").Value = 18513.143
$ws.Range("M32").Value = -5524.5757
$ws.Range("N32").Value = -19087.143
$ws.Range("H37").Value = 10759.333
$ws.Range("I37").Value = 4889
$ws.Range("K37").Value = 4889
$ws.Range("M37").Value = -4616
$ws.Range("H61").Value = 2231.5
$ws.Range("I61").Value = 2231.5
$ws.Range("K61").Value = 2231.5
$ws.Range("M61").Value = -2019.5
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H97").Value = 861.375
$ws.Range("I97").Value = 846.7692
$ws.Range("K97").Value = 846.7692
$ws.Range("M97").Value = -350.7692
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 1003.3333
$ws.Range("I122").Value = 1003.3333
$ws.Range("K122").Value = 3009.9999
$ws.Range("M122").Value = -559.9998999999998
$ws.Range("H136").Value = 2231.5
$ws.Range("I136").Value = 2231.5
$ws.Range("K136").Value = 6694.5
$ws.Range("M136").Value = -4144.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 3702
$ws.Range("J100").Value = 3702
$ws.Range("L100").Value = 3702
$ws.Range("N100").Value = -5866
$ws.Range("H105").Value = 1849.625
$ws.Range("I105").Value = 1608.7273
$ws.Range("K105").Value = 1608.7273
$ws.Range("M105").Value = 138.2727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2846
$ws.Range("I31").Value = 1463.4
$ws.Range("J31").Value = 6302.5
$ws.Range("K31").Value = 1463.4
$ws.Range("L31").Value = 6302.5
$ws.Range("M31").Value = -1168.4
$ws.Range("N31").Value = -6892.5
$ws.Range("H34").Value = 2846
$ws.Range("I34").Value = 1463.4
$ws.Range("J34").Value = 6302.5
$ws.Range("K34").Value = 1463.4
$ws.Range("L34").Value = 6302.5
$ws.Range("M34").Value = -1261.4
$ws.Range("N34").Value = -6706.5
$ws.Range("H41").Value = 26308.334
$ws.Range("I41").Value = 10750
$ws.Range("J41").Value = 29420
$ws.Range("K41").Value = 10750
$ws.Range("L41").Value = 29420
$ws.Range("M41").Value = -10322
$ws.Range("N41").Value = -30276
$ws.Range("H53").Value = 63999
$ws.Range("J53").Value = 63999
$ws.Range("L53").Value = 63999
$ws.Range("N53").Value = -65213
$ws.Range("H105").Value = 3833.147
$ws.Range("I105").Value = 657.1667
$ws.Range("J105").Value = 7406.125
$ws.Range("K105").Value = 657.1667
$ws.Range("L105").Value = 7406.125
$ws.Range("M105").Value = 1089.8333
$ws.Range("N105").Value = -10900.125
$ws.Range("H106").Value = 7777
$ws.Range("J106").Value = 7777
$ws.Range("L106").Value = 7777
$ws.Range("N106").Value = -10301
$ws.Range("H141").Value = 134584.28
$ws.Range("J141").Value = 134584.28
$ws.Range("L141").Value = 134584.28
$ws.Range("N141").Value = -144944.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1501362.6
$ws.Range("I4").Value = 2625593.5
$ws.Range("J4").Value = 2388
$ws.Range("K4").Value = 7876780.5
$ws.Range("L4").Value = 7164
$ws.Range("M4").Value = -7876668.5
$ws.Range("N4").Value = -7388
$ws.Range("H33").Value = 49.416668
$ws.Range("I33").Value = 22.833334
$ws.Range("J33").Value = 76
$ws.Range("K33").Value = 137.000004
$ws.Range("L33").Value = 456
$ws.Range("M33").Value = 145.999996
$ws.Range("N33").Value = -1022
$ws.Range("H131").Value = 1073.6222
$ws.Range("I131").Value = 800
$ws.Range("J131").Value = 1086.3489
$ws.Range("K131").Value = 2400
$ws.Range("L131").Value = 3259.0467
$ws.Range("M131").Value = 2640
$ws.Range("N131").Value = -13339.0467

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 234.6
$ws.Range("I31").Value = 234.6
$ws.Range("K31").Value = 234.6
$ws.Range("M31").Value = 57.40000000000001
$ws.Range("H37").Value = 234.6
$ws.Range("I37").Value = 234.6
$ws.Range("K37").Value = 234.6
$ws.Range("M37").Value = 42.40000000000001
$ws.Range("H98").Value = 7296.3335
$ws.Range("J98").Value = 7296.3335
$ws.Range("L98").Value = 7296.3335
$ws.Range("N98").Value = -13286.3335
$ws.Range("H113").Value = 3528.0588
$ws.Range("I113").Value = 2623.375
$ws.Range("K113").Value = 2623.375
$ws.Range("M113").Value = -453.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1634
$ws.Range("H132").Value = 2378.2666
$ws.Range("I132").Value = 1869.9
$ws.Range("K132").Value = 5609.700000000001
$ws.Range("M132").Value = -3079.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H105").Value = 47807.5
$ws.Range("J105").Value = 47807.5
$ws.Range("L105").Value = 47807.5
$ws.Range("N105").Value = -54795.5
$ws.Range("H122").Value = 2212.9
$ws.Range("I122").Value = 2212.9
$ws.Range("K122").Value = 6638.700000000001
$ws.Range("M122").Value = -4188.700000000001
$ws.Range("H136").Value = 1474.1111
$ws.Range("I136").Value = 968.93335
$ws.Range("K136").Value = 2906.80005
$ws.Range("M136").Value = -356.8000499999998
